# Applies the edits described by the commit "PROJ: xslx updated, readme
# initial version, simplified config" to the autoavaliacao workbook.
#
# Summary of the substantive changes on the "Autoavaliação" sheet:
#   - B4 now holds the student/group identifier "RCL12"
#   - B5 now holds the total grade 86
#   - The per-criterion self-assessment cells in columns B and F (rows
#     9-26) are filled in with the same values the grader already put in
#     the neighbouring C/G "max" columns (so the SUM()-based totals in
#     B7/B8/F8 recompute automatically)
#   - The checklist in column B (rows 32-55) is marked "Completo"
#   - The active selection on the sheet moves to H22

$wb = $excel.ActiveWorkbook
# First sheet is "Autoavaliação"; index access sidesteps any accented-name
# encoding pitfalls in the host shell.
$ws = $wb.Worksheets.Item(1)
$ws.Select()

# --- Header block -----------------------------------------------------
$ws.Range("B4").Value = "RCL12"
$ws.Range("B5").Value = 86

# --- Left table (columns A-C): self-assessment scores ------------------
$leftScores = @{
    9  = 1
    10 = 0.5
    11 = 0.5
    12 = 1
    13 = 1
    14 = 0.25
    15 = 0.25
    17 = 1.25
    21 = 2
    22 = 0.5
    23 = 1
    24 = 0.5
    26 = 2
}
foreach ($row in $leftScores.Keys) {
    $ws.Cells.Item($row, 2).Value = $leftScores[$row]
}

# --- Right table (columns E-G): self-assessment scores ------------------
$rightScores = @{
    9  = 0.25
    10 = 0.5
    11 = 0.5
    12 = 0.75
    13 = 1
    14 = 0.25
    15 = 0.25
    17 = 1
    21 = 0.75
    22 = 0.5
    23 = 0.75
    24 = 0.5
    26 = 0.75
}
foreach ($row in $rightScores.Keys) {
    $ws.Cells.Item($row, 6).Value = $rightScores[$row]
}

# --- AS server test checklist: mark every item "Completo" --------------
for ($row = 32; $row -le 55; $row++) {
    $ws.Cells.Item($row, 2).Value = "Completo"
}

# --- Move the selection, matching the saved view in the edited file ----
$ws.Range("H22").Select()
